$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F7").Value = 777
$ws1.Range("G7").Value = 45
$ws1.Range("G8").Value = 60
$ws1.Range("F9").Value = 200
$ws1.Range("G9").Value = 70
$ws1.Range("F11").Value = 478
$ws1.Range("F13").Value = 134
$ws1.Range("F14").Value = 130
$ws1.Range("F17").Value = 103
$ws1.Range("F18").Value = 683
$ws1.Range("F21").Value = 273
$ws1.Range("F23").Value = 6162
$ws1.Range("F25").Value = 131
$ws1.Range("F26").Value = 129
$ws1.Range("F28").Value = 14915
$ws1.Range("F29").Value = 1482
$ws1.Range("F33").Value = 10872
$ws1.Range("F34").Value = 689
$ws1.Range("F35").Value = 4249
$ws1.Range("F36").Value = 190
$ws1.Range("F38").Value = 118

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 777
$ws4.Range("G7").Value = 45
$ws4.Range("G8").Value = 60
$ws4.Range("F9").Value = 200
$ws4.Range("G9").Value = 70
$ws4.Range("F11").Value = 478
$ws4.Range("F13").Value = 134
$ws4.Range("F14").Value = 130
$ws4.Range("F18").Value = 103
$ws4.Range("F19").Value = 683
$ws4.Range("F23").Value = 273
$ws4.Range("F24").Value = 0
$ws4.Range("F26").Value = 6162
$ws4.Range("F28").Value = 131
$ws4.Range("F29").Value = 129
$ws4.Range("F31").Value = 14915
$ws4.Range("F32").Value = 1482
$ws4.Range("F36").Value = 10872
$ws4.Range("F37").Value = 689
$ws4.Range("F38").Value = 4249
$ws4.Range("F39").Value = 190
$ws4.Range("F41").Value = 118
